$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo in D7: "abhijit291@gmail.com" -> "abhiit291@gmail.com"
$ws.Range("D7").Value = "abhiit291@gmail.com"

# Make D7 the active selection (matches author's last edit position)
$ws.Range("D7").Select()
